# Applies scheduled-runner market-price/profit updates to the Sagittarius_Profits workbook.
# Each class sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) stores hardcoded (non-formula)
# market snapshot values in columns H:N; this script overwrites the specific cells
# that changed between runs, matching the upstream commit exactly.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1541.4286
$ws.Range("I33").Value = 1723.3334
$ws.Range("K33").Value = 1723.3334
$ws.Range("M33").Value = -1494.3334
$ws.Range("H41").Value = 1010
$ws.Range("I41").Value = 20
$ws.Range("K41").Value = 20
$ws.Range("M41").Value = 420
$ws.Range("H98").Value = 5154.609
$ws.Range("I98").Value = 993.7143
$ws.Range("K98").Value = 993.7143
$ws.Range("M98").Value = 504.2857
$ws.Range("H99").Value = 349
$ws.Range("I99").Value = 349
$ws.Range("K99").Value = 1047
$ws.Range("M99").Value = 451
$ws.Range("H112").Value = 1099.0555
$ws.Range("J112").Value = 1075.4706
$ws.Range("L112").Value = 3226.4118
$ws.Range("N112").Value = -5442.4118
$ws.Range("H122").Value = 5154.609
$ws.Range("I122").Value = 993.7143
$ws.Range("K122").Value = 2981.1429
$ws.Range("M122").Value = -531.1428999999998
$ws.Range("H131").Value = 11174.9375
$ws.Range("I131").Value = 2467.3333
$ws.Range("J131").Value = 16399.5
$ws.Range("K131").Value = 7401.999899999999
$ws.Range("L131").Value = 49198.5
$ws.Range("M131").Value = -2361.999899999999
$ws.Range("N131").Value = -59278.5
$ws.Range("H138").Value = 3681.9368
$ws.Range("I138").Value = 2681.1292
$ws.Range("J138").Value = 4166.703
$ws.Range("K138").Value = 8043.3876
$ws.Range("L138").Value = 12500.109
$ws.Range("M138").Value = -2903.3876
$ws.Range("N138").Value = -22780.109

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4223.162
$ws.Range("I32").Value = 3721.8
$ws.Range("J32").Value = 12997
$ws.Range("K32").Value = 3721.8
$ws.Range("L32").Value = 12997
$ws.Range("M32").Value = -3434.8
$ws.Range("N32").Value = -13571
$ws.Range("H97").Value = 726
$ws.Range("J97").Value = 418.66666
$ws.Range("L97").Value = 418.66666
$ws.Range("N97").Value = -1410.66666

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 328.8
$ws.Range("I94").Value = 320.8889
$ws.Range("K94").Value = 320.8889
$ws.Range("M94").Value = 130.1111
$ws.Range("H96").Value = 10479.625
$ws.Range("I96").Value = 10479.625
$ws.Range("K96").Value = 10479.625
$ws.Range("M96").Value = -7733.625
$ws.Range("H99").Value = 3179
$ws.Range("I99").Value = 2565.3076
$ws.Range("K99").Value = 2565.3076
$ws.Range("M99").Value = -1067.3076
$ws.Range("H105").Value = 3222
$ws.Range("I105").Value = 3298
$ws.Range("J105").Value = 2994
$ws.Range("K105").Value = 3298
$ws.Range("L105").Value = 2994
$ws.Range("M105").Value = -1551
$ws.Range("N105").Value = -6488
$ws.Range("H107").Value = 3436.25
$ws.Range("I107").Value = 2615.6
$ws.Range("K107").Value = 2615.6
$ws.Range("M107").Value = -695.5999999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7123.7856
$ws.Range("I31").Value = 6024
$ws.Range("K31").Value = 6024
$ws.Range("M31").Value = -5729
$ws.Range("H34").Value = 7123.7856
$ws.Range("I34").Value = 6024
$ws.Range("K34").Value = 6024
$ws.Range("M34").Value = -5822
$ws.Range("H105").Value = 2164.0527
$ws.Range("I105").Value = 2219.3845
$ws.Range("J105").Value = 2044.1666
$ws.Range("K105").Value = 2219.3845
$ws.Range("L105").Value = 2044.1666
$ws.Range("M105").Value = -472.3845000000001
$ws.Range("N105").Value = -5538.1666

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1050.5
$ws.Range("J59").Value = 1466.3334
$ws.Range("L59").Value = 4399.0002
$ws.Range("N59").Value = -5479.0002
$ws.Range("H86").Value = 7766.3335
$ws.Range("J86").Value = 11500
$ws.Range("L86").Value = 34500
$ws.Range("N86").Value = -36872
$ws.Range("H89").Value = 7766.3335
$ws.Range("J89").Value = 11500
$ws.Range("L89").Value = 103500
$ws.Range("N89").Value = -115356
$ws.Range("H98").Value = 664.2353000000001
$ws.Range("J98").Value = 651.4
$ws.Range("L98").Value = 1954.2
$ws.Range("N98").Value = -4950.2
$ws.Range("H121").Value = 12560.5625
$ws.Range("I121").Value = 18686.857
$ws.Range("K121").Value = 56060.571
$ws.Range("M121").Value = -54750.571
$ws.Range("H130").Value = 2506.6667
$ws.Range("I130").Value = 2015
$ws.Range("K130").Value = 6045
$ws.Range("M130").Value = -1025

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 59999
$ws.Range("I33").Value = 59999
$ws.Range("K33").Value = 59999
$ws.Range("M33").Value = -59747
$ws.Range("H70").Value = 5766.2856
$ws.Range("I70").Value = 5488
$ws.Range("K70").Value = 5488
$ws.Range("M70").Value = -5218
$ws.Range("H73").Value = 5766.2856
$ws.Range("I73").Value = 5488
$ws.Range("K73").Value = 5488
$ws.Range("M73").Value = -4552
$ws.Range("H80").Value = 10033.286
$ws.Range("I80").Value = 2527.4285
$ws.Range("K80").Value = 2527.4285
$ws.Range("M80").Value = -1529.4285
$ws.Range("H83").Value = 10033.286
$ws.Range("I83").Value = 2527.4285
$ws.Range("K83").Value = 12637.1425
$ws.Range("M83").Value = -7645.1425
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1098.8334
$ws.Range("I16").Value = 923.5
$ws.Range("J16").Value = 1449.5
$ws.Range("K16").Value = 923.5
$ws.Range("L16").Value = 1449.5
$ws.Range("M16").Value = -753.5
$ws.Range("N16").Value = -1789.5
$ws.Range("H55").Value = 1106.4762
$ws.Range("I55").Value = 1013.3333
$ws.Range("J55").Value = 1230.6666
$ws.Range("K55").Value = 1013.3333
$ws.Range("L55").Value = 1230.6666
$ws.Range("M55").Value = -840.3333
$ws.Range("N55").Value = -1576.6666
$ws.Range("H61").Value = 6570.857
$ws.Range("I61").Value = 5582.9443
$ws.Range("J61").Value = 12498.333
$ws.Range("K61").Value = 5582.9443
$ws.Range("L61").Value = 12498.333
$ws.Range("M61").Value = -5380.9443
$ws.Range("N61").Value = -12902.333
$ws.Range("H68").Value = 2699.3125
$ws.Range("I68").Value = 1927.1111
$ws.Range("K68").Value = 1927.1111
$ws.Range("M68").Value = -1178.1111
$ws.Range("H71").Value = 2699.3125
$ws.Range("I71").Value = 1927.1111
$ws.Range("K71").Value = 9635.5555
$ws.Range("M71").Value = -5891.5555
$ws.Range("H100").Value = 4700
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959
$ws.Range("H113").Value = 6570.857
$ws.Range("I113").Value = 5582.9443
$ws.Range("J113").Value = 12498.333
$ws.Range("K113").Value = 5582.9443
$ws.Range("L113").Value = 12498.333
$ws.Range("M113").Value = -3412.9443
$ws.Range("N113").Value = -16838.333
$ws.Range("H132").Value = 3899
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7548.4
$ws.Range("I132").Value = 7497.143
$ws.Range("J132").Value = 7668
$ws.Range("K132").Value = 22491.429
$ws.Range("L132").Value = 23004
$ws.Range("M132").Value = -19961.429
$ws.Range("N132").Value = -28064
